$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.236.96'
$ws.Range("E2").Value = '  +0.64%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.603.06'
$ws.Range("E3").Value = '  -0.04%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9999'
$ws.Range("E5").Value = '  -0.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.37'
$ws.Range("E6").Value = '  +0.64%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3773'
$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3634'
$ws.Range("E8").Value = '  -0.58%  '

$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.49'
$ws.Range("E9").Value = '  +3.01%  '

$ws.Range("E10").Value = '  +0.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08140'
$ws.Range("E11").Value = '  -0.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9998'
$ws.Range("E12").Value = '  -0.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.91'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.602'
$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.408'
$ws.Range("E15").Value = '  +0.23%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001249'
$ws.Range("E16").Value = '  -0.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.600.76'
$ws.Range("E17").Value = '  -0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.88'
$ws.Range("E18").Value = '  +1.86%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06893'
$ws.Range("E19").Value = '  +0.07%  '

$ws.Range("E20").Value = '  -0.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.533'
$ws.Range("E21").Value = '  -0.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  -0.22%  '

$ws.Range("E23").Value = '  -1.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.230.58'
$ws.Range("E24").Value = '  +0.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.389'
$ws.Range("E25").Value = '  +0.98%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.999'
$ws.Range("E26").Value = '  +7.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.22'
$ws.Range("E27").Value = '  +0.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.10'
$ws.Range("E28").Value = '  -0.24%  '

$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.18'
$ws.Range("E30").Value = '  -0.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.382'
$ws.Range("E31").Value = '  +1.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.767'
$ws.Range("E32").Value = '  -1.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.780.53'
$ws.Range("E33").Value = '  +0.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9654'
$ws.Range("E34").Value = '  +0.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07536'
$ws.Range("E35").Value = '  -1.99%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02731'
$ws.Range("E36").Value = '  +0.44%  '

$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.26'
$ws.Range("E37").Value = '  -2.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2532'
$ws.Range("E38").Value = '  -0.64%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.08805'
$ws.Range("E39").Value = '  -1.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.096'
$ws.Range("E40").Value = '  -3.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.397'
$ws.Range("E41").Value = '  +1.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7108'
$ws.Range("E42").Value = '  +0.34%  '

$ws.Range("E43").Value = '  -1.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.55'
$ws.Range("E44").Value = '  +0.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6546'
$ws.Range("E45").Value = '  -1.48%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.323'
$ws.Range("E46").Value = '  -0.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.014'
$ws.Range("E47").Value = '  +0.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.56'
$ws.Range("E48").Value = '  +0.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07968'
$ws.Range("E49").Value = '  +0.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.206'
$ws.Range("E50").Value = '  -3.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.206'
$ws.Range("E51").Value = '  -0.22%  '
